# Apply the "fixed 2 sheet code, found a bug in escape" edit.
#
# Summary of changes:
#  1. Workbook window position (xWindow/yWindow) moved.
#  2. Sheet "PartI": the "related_mats" column (L) is deleted, shifting
#     everything after it one column to the left.
#  3. Three shared strings are renamed:
#       collection_date -> "collection date"
#       acq_info        -> "acqinfo"
#       c_level         -> "c0x level"
#  4. Sheet "PartII": column D gets a Text ("@") number format (matching
#     the existing column F formatting), and the selections on both sheets
#     move.

$wb = $excel.ActiveWorkbook

# --- 1. Workbook window position -------------------------------------
try {
    $win = $wb.Windows.Item(1)
    $win.Left = 2520
    $win.Top = 1640
} catch {
    # Not fatal if the host doesn't expose window placement.
}

# --- 2 & 3. Sheet "PartI" ----------------------------------------------
$ws1 = $wb.Worksheets.Item("PartI")

# Delete column L ("related_mats") entirely; remaining columns shift left.
$ws1.Columns.Item(12).Delete()

# Rename the shared strings used on this sheet.
$ws1.Cells.Item(1, 3).Value = "collection date"   # was "collection_date"
$ws1.Cells.Item(1, 9).Value = "acqinfo"           # was "acq_info"

# Move the active selection to L1.
$ws1.Activate()
$ws1.Range("L1").Select()

# --- 4. Sheet "PartII" --------------------------------------------------
$ws2 = $wb.Worksheets.Item("PartII")

# Rename the "c_level" shared string.
$ws2.Cells.Item(1, 2).Value = "c0x level"

# Column D now gets the same Text number format as column F.
$ws2.Columns.Item(4).NumberFormat = "@"

# Move the active selection to D2.
$ws2.Activate()
$ws2.Range("D2").Select()
